$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Auftragsfenster")

$ws.Range("B5").Value = $ws.Range("B4").Value()
$ws.Range("C5").Value = $ws.Range("C4").Value()
$ws.Range("D5").Value = $ws.Range("D4").Value()
$ws.Range("E5").Value = "Programmvorschau2022"
$ws.Range("G5").Value = $ws.Range("G4").Value()
$ws.Range("H5").Value = $ws.Range("H4").Value()
$ws.Range("I5").Value = $ws.Range("I4").Value()
$ws.Range("L5").Value = $ws.Range("L4").Value()

$ws.Range("B6").Value = $ws.Range("B4").Value()
$ws.Range("C6").Value = $ws.Range("C4").Value()
$ws.Range("D6").Value = $ws.Range("D4").Value()
$ws.Range("E6").Value = "Trailer2022"
$ws.Range("G6").Value = $ws.Range("G4").Value()
$ws.Range("H6").Value = $ws.Range("H4").Value()
$ws.Range("I6").Value = $ws.Range("I4").Value()
$ws.Range("L6").Value = $ws.Range("L4").Value()
